$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new wishlist entry as the next row (row 20).
$ws.Range("A20").Value = "Diez lunas para una espera"
$ws.Range("B20").Value = "Velia Vidal"

# The Editorial column is left blank for this entry (matches the existing
# empty-Editorial rows above it), but we still want the cell itself to exist
# in the sheet rather than be entirely absent. Setting the style (without
# changing it) is enough to materialize an empty cell at C20, the same way
# the other blank Editorial cells are materialized as empty cells.
$ws.Range("C20").Style = "Normal"
